# Regenerate the handback status report for the new handoff batch:
#   072e5562-1f19-4667-b4f4-49482e7ffe8d -> 71862875-3cce-4686-b8aa-acc814e9dcde
#   dbff55ea-36b5-4ee3-b886-a544221a78ad -> ffff47e1dfb2-2cbc-4270-b719-7b65a80e6f18
# and refreshed handoff/handback timestamps + xlf hashes.

$wb = $excel.ActiveWorkbook

$newUuid1 = "71862875-3cce-4686-b8aa-acc814e9dcde"
$newUuid2 = "ffff47e1dfb2-2cbc-4270-b719-7b65a80e6f18"
$newHash  = "6498da75e234a9e78b2721ed45546758b875fb12"

$file1Name = "$newUuid1.md"
$file2Name = "$newUuid2.md"
$file1Path = "e2e\$newUuid1.md"
$file2Path = "e2e\$newUuid2.md"

$zhXlf = "$newUuid1.$newHash.zh-cn.xlf"
$deXlf = "$newUuid1.$newHash.de-de.xlf"

$genDate      = "2016-08-15 22:58:29"
$zhHandoffDt  = "2016-08-15 22:58:24"
$zhHandbackDt = "2016-08-15 22:58:41"
$deHandbackDt = "2016-08-15 22:58:48"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $file1Name
$ws1.Range("G2").Value = $genDate
$ws1.Range("A3").Value = $file2Name
$ws1.Range("G3").Value = $genDate

$h = $ws1.Range("B2").Hyperlinks.Item(1)
$h.TextToDisplay = $file1Path
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acb16955c033ccde0bbec3579334f9e43841a156/e2e/072e5562-1f19-4667-b4f4-49482e7ffe8d.md"

$h = $ws1.Range("B3").Hyperlinks.Item(1)
$h.TextToDisplay = $file2Path
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acb16955c033ccde0bbec3579334f9e43841a156/e2e/dbff55ea-36b5-4ee3-b886-a544221a78ad.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $file1Name
$ws2.Range("I2").Value = $file1Name
$ws2.Range("G2").Value = $zhXlf
$ws2.Range("J2").Value = $zhXlf
$ws2.Range("H2").Value = $zhHandoffDt
$ws2.Range("K2").Value = $zhHandbackDt

$ws2.Range("A3").Value = $file2Name
$ws2.Range("I3").Value = $file2Name
$ws2.Range("G3").Value = $zhXlf
$ws2.Range("J3").Value = $zhXlf
$ws2.Range("H3").Value = $zhHandoffDt
$ws2.Range("K3").Value = $zhHandbackDt

$h = $ws2.Range("A2").Hyperlinks.Item(1)
$h.TextToDisplay = $file1Name
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acb16955c033ccde0bbec3579334f9e43841a156/e2e/072e5562-1f19-4667-b4f4-49482e7ffe8d.md"

$h = $ws2.Range("I2").Hyperlinks.Item(1)
$h.TextToDisplay = $file1Name
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bc4961a2adaefc8d3fa2fdb89e58137f5dc825ce/e2e/072e5562-1f19-4667-b4f4-49482e7ffe8d.md"

$h = $ws2.Range("A3").Hyperlinks.Item(1)
$h.TextToDisplay = $file2Name
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acb16955c033ccde0bbec3579334f9e43841a156/e2e/dbff55ea-36b5-4ee3-b886-a544221a78ad.md"

$h = $ws2.Range("I3").Hyperlinks.Item(1)
$h.TextToDisplay = $file2Name
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bc4961a2adaefc8d3fa2fdb89e58137f5dc825ce/e2e/dbff55ea-36b5-4ee3-b886-a544221a78ad.md"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $file1Name
$ws3.Range("I2").Value = $file1Name
$ws3.Range("G2").Value = $deXlf
$ws3.Range("J2").Value = $deXlf
$ws3.Range("H2").Value = $genDate
$ws3.Range("K2").Value = $deHandbackDt

$ws3.Range("A3").Value = $file2Name
$ws3.Range("I3").Value = $file2Name
$ws3.Range("G3").Value = $deXlf
$ws3.Range("J3").Value = $deXlf
$ws3.Range("H3").Value = $genDate
$ws3.Range("K3").Value = $deHandbackDt

$h = $ws3.Range("A2").Hyperlinks.Item(1)
$h.TextToDisplay = $file1Name
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acb16955c033ccde0bbec3579334f9e43841a156/e2e/072e5562-1f19-4667-b4f4-49482e7ffe8d.md"

$h = $ws3.Range("I2").Hyperlinks.Item(1)
$h.TextToDisplay = $file1Name
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8c4d8383d381dbecc0463eb1db02c85d0ec7fef0/e2e/072e5562-1f19-4667-b4f4-49482e7ffe8d.md"

$h = $ws3.Range("A3").Hyperlinks.Item(1)
$h.TextToDisplay = $file2Name
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acb16955c033ccde0bbec3579334f9e43841a156/e2e/dbff55ea-36b5-4ee3-b886-a544221a78ad.md"

$h = $ws3.Range("I3").Hyperlinks.Item(1)
$h.TextToDisplay = $file2Name
$h.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8c4d8383d381dbecc0463eb1db02c85d0ec7fef0/e2e/dbff55ea-36b5-4ee3-b886-a544221a78ad.md"
